$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-31 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-01 Monday", 2) | Out-Null
$d.Content.Find.Execute("84×67=", $true, $false, $false, $false, $false, $true, 1, $false, "34×17=", 2) | Out-Null
$d.Content.Find.Execute("39×97=", $true, $false, $false, $false, $false, $true, 1, $false, "91×39=", 2) | Out-Null
$d.Content.Find.Execute("17×95=", $true, $false, $false, $false, $false, $true, 1, $false, "40×68=", 2) | Out-Null
$d.Content.Find.Execute("94×29=", $true, $false, $false, $false, $false, $true, 1, $false, "93×27=", 2) | Out-Null
$d.Content.Find.Execute("42×71=", $true, $false, $false, $false, $false, $true, 1, $false, "47×19=", 2) | Out-Null
$d.Content.Find.Execute("88×21=", $true, $false, $false, $false, $false, $true, 1, $false, "18×24=", 2) | Out-Null
$d.Content.Find.Execute("47×16=", $true, $false, $false, $false, $false, $true, 1, $false, "46×53=", 2) | Out-Null
$d.Content.Find.Execute("31×77=", $true, $false, $false, $false, $false, $true, 1, $false, "68×16=", 2) | Out-Null
$d.Content.Find.Execute("14×87=", $true, $false, $false, $false, $false, $true, 1, $false, "16×79=", 2) | Out-Null
$d.Content.Find.Execute("65×11=", $true, $false, $false, $false, $false, $true, 1, $false, "71×41=", 2) | Out-Null
$d.Content.Find.Execute("90×91=", $true, $false, $false, $false, $false, $true, 1, $false, "61×31=", 2) | Out-Null
$d.Content.Find.Execute("68×58=", $true, $false, $false, $false, $false, $true, 1, $false, "85×35=", 2) | Out-Null
$d.Content.Find.Execute("83×92=", $true, $false, $false, $false, $false, $true, 1, $false, "70×23=", 2) | Out-Null
$d.Content.Find.Execute("37×65=", $true, $false, $false, $false, $false, $true, 1, $false, "20×95=", 2) | Out-Null
$d.Content.Find.Execute("57×39=", $true, $false, $false, $false, $false, $true, 1, $false, "82×82=", 2) | Out-Null
$d.Content.Find.Execute("11×53=", $true, $false, $false, $false, $false, $true, 1, $false, "73×34=", 2) | Out-Null
$d.Content.Find.Execute("54×80=", $true, $false, $false, $false, $false, $true, 1, $false, "38×78=", 2) | Out-Null
$d.Content.Find.Execute("22×14=", $true, $false, $false, $false, $false, $true, 1, $false, "87×88=", 2) | Out-Null
$d.Content.Find.Execute("11×36=", $true, $false, $false, $false, $false, $true, 1, $false, "81×76=", 2) | Out-Null
$d.Content.Find.Execute("56×38=", $true, $false, $false, $false, $false, $true, 1, $false, "41×47=", 2) | Out-Null
$d.Content.Find.Execute("62×52=", $true, $false, $false, $false, $false, $true, 1, $false, "72×39=", 2) | Out-Null
$d.Content.Find.Execute("29×86=", $true, $false, $false, $false, $false, $true, 1, $false, "60×48=", 2) | Out-Null
$d.Content.Find.Execute("90×55=", $true, $false, $false, $false, $false, $true, 1, $false, "60×49=", 2) | Out-Null
$d.Content.Find.Execute("46×19=", $true, $false, $false, $false, $false, $true, 1, $false, "96×11=", 2) | Out-Null
$d.Content.Find.Execute("92×50=", $true, $false, $false, $false, $false, $true, 1, $false, "66×25=", 2) | Out-Null
